$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Borrado múltiple": delete the empty leading columns (A:E) and leading rows (1:3)
# so that the "cuentas" table shifts from F4:H17 up to A1:C14.
$ws.Columns("A:E").Delete()
$ws.Rows("1:3").Delete()

# Update the "Cuentas" defined name to track the table's new location.
$wb.Names("Cuentas").RefersTo = "=cuentas!`$A`$1:`$B`$8"

# Leave the selection where the user ended up after the deletion.
[void]$ws.Range("B16").Select()
